# Update the "Latest HO Xliff Generate Date", "Correspond Handoff Datetime"
# and "Correspond Handback DateTime" timestamps produced by a fresh handback
# report generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for b095d723... row
$wsOverview.Range("G4").Value = "2016-09-03 18:58:18"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for b095d723... row
$wsZhCn.Range("H4").Value = "2016-09-03 18:58:14"
$wsZhCn.Range("K4").Value = "2016-09-03 18:58:31"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime for b095d723... row
$wsDeDe.Range("H4").Value = "2016-09-03 18:58:18"
$wsDeDe.Range("K4").Value = "2016-09-03 18:58:38"
